# The pptx writer now uses the incoming/explicit column widths for a
# table instead of always dividing the available width evenly between
# columns. For the first table on slide 6 this means both columns grow
# from 197pt (2501900 EMU) to 198pt (2514600 EMU) each - matching the
# width already used by the second (already-fixed) table on the same
# slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

# Shape 3 is the first "Content Placeholder" graphicFrame (the table
# positioned at x=1117600) that still has the old, unevenly-divided
# column widths.
$shape = $s.Shapes.Item(3)

if ($shape.HasTable) {
    $table = $shape.Table
    for ($c = 1; $c -le $table.Columns.Count; $c++) {
        $table.Columns.Item($c).Width = 198
    }
}
